$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4175
$ws.Range("E2").Value = 145
$ws.Range("F2").Value = 145
$ws.Range("G2").Value = 95
$ws.Range("H2").Value = 119
$ws.Range("I2").Value = 119
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 6203
$ws.Range("L2").Value = 2874
$ws.Range("M2").Value = 3330
$ws.Range("N2").Value = 3316
$ws.Range("O2").Value = 14
$ws.Range("P2").Value = 251
$ws.Range("Q2").Value = 394
$ws.Range("R2").Value = -459
$ws.Range("S2").Value = -256
$ws.Range("T2").Value = 303
$ws.Range("U2").Value = 90
$ws.Range("V2").Value = 1443
$ws.Range("W2").Value = 3.46
$ws.Range("X2").Value = 2.86
$ws.Range("Y2").Value = 3.61
$ws.Range("Z2").Value = 1.91
$ws.Range("AA2").Value = 86.32
$ws.Range("AB2").Value = 1228.3
$ws.Range("AC2").Value = 451
$ws.Range("AD2").Value = 31.43
$ws.Range("AE2").Value = 13023
$ws.Range("AF2").Value = 1.09
$ws.Range("AG2").Value = 191
$ws.Range("AH2").Value = 1.35
$ws.Range("AI2").Value = 40.89
$ws.Range("AJ2").Value = 26292258

# Row 3
$ws.Range("D3").Value = 383
$ws.Range("E3").Value = -49
$ws.Range("F3").Value = 237
$ws.Range("G3").Value = -45
$ws.Range("H3").Value = 213
$ws.Range("I3").Value = 210
$ws.Range("J3").Value = 4
$ws.Range("K3").Value = 7056
$ws.Range("L3").Value = 3557
$ws.Range("M3").Value = 3500
$ws.Range("N3").Value = 3481
$ws.Range("O3").Value = 19
$ws.Range("P3").Value = 251
$ws.Range("Q3").Value = 561
$ws.Range("R3").Value = -399
$ws.Range("S3").Value = 420
$ws.Range("T3").Value = 152
$ws.Range("U3").Value = 409
$ws.Range("V3").Value = 1915
$ws.Range("W3").Value = -12.89
$ws.Range("X3").Value = 55.69
$ws.Range("Y3").Value = 6.17
$ws.Range("Z3").Value = 3.22
$ws.Range("AA3").Value = 101.62
$ws.Range("AB3").Value = 1284.29
$ws.Range("AC3").Value = 798
$ws.Range("AD3").Value = 29.74
$ws.Range("AE3").Value = 13672
$ws.Range("AF3").Value = 1.74
$ws.Range("AG3").Value = 286
$ws.Range("AH3").Value = 1.21
$ws.Range("AI3").Value = 34.66
$ws.Range("AJ3").Value = 26292258

# Row 4
$ws.Range("D4").Value = 408
$ws.Range("E4").Value = -37
$ws.Range("F4").Value = -35
$ws.Range("G4").Value = -46
$ws.Range("H4").Value = 2662
$ws.Range("I4").Value = 2660
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 2041
$ws.Range("L4").Value = 972
$ws.Range("M4").Value = 1068
$ws.Range("N4").Value = 1048
$ws.Range("O4").Value = 21
$ws.Range("P4").Value = 72
$ws.Range("Q4").Value = 466
$ws.Range("R4").Value = -424
$ws.Range("S4").Value = -195
$ws.Range("T4").Value = 284
$ws.Range("U4").Value = 182
$ws.Range("V4").Value = 636
$ws.Range("W4").Value = -9.16
$ws.Range("X4").Value = 652.89
$ws.Range("Y4").Value = 117.47
$ws.Range("Z4").Value = 58.52
$ws.Range("AA4").Value = 91.03
$ws.Range("AB4").Value = 8151.55
$ws.Range("AC4").Value = 14323
$ws.Range("AD4").Value = 1.72
$ws.Range("AE4").Value = 14292
$ws.Range("AF4").Value = 1.73
$ws.Range("AG4").Value = 143
$ws.Range("AH4").Value = 0.58
$ws.Range("AI4").Value = 0.39
$ws.Range("AJ4").Value = 7575536

# Row 5
$ws.Range("D5").Value = 411
$ws.Range("E5").Value = -18
$ws.Range("F5").Value = 7
$ws.Range("G5").Value = -8
$ws.Range("H5").Value = -6
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = -6
$ws.Range("K5").Value = 2742
$ws.Range("L5").Value = 1002
$ws.Range("M5").Value = 1741
$ws.Range("N5").Value = 1727
$ws.Range("O5").Value = 14
$ws.Range("P5").Value = 105
$ws.Range("Q5").Value = -22
$ws.Range("R5").Value = -191
$ws.Range("S5").Value = 9
$ws.Range("T5").Value = 21
$ws.Range("U5").Value = -43
$ws.Range("V5").Value = 660
$ws.Range("W5").Value = -4.45
$ws.Range("X5").Value = -1.35
$ws.Range("Y5").Value = 0.05
$ws.Range("Z5").Value = -0.23
$ws.Range("AA5").Value = 57.55
$ws.Range("AB5").Value = 1541.59
$ws.Range("AC5").Value = 7
$ws.Range("AD5").Value = 1998.32
$ws.Range("AE5").Value = 16049
$ws.Range("AF5").Value = 0.9
$ws.Range("AG5").Value = 191
$ws.Range("AH5").Value = 1.32
$ws.Range("AI5").Value = 2775.88
$ws.Range("AJ5").Value = 11002587

# Row 6
$ws.Range("D6").Value = 470
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = -7
$ws.Range("G6").Value = -29
$ws.Range("H6").Value = -38
$ws.Range("I6").Value = -37
$ws.Range("K6").Value = 2655
$ws.Range("L6").Value = 1005
$ws.Range("M6").Value = 1649
$ws.Range("N6").Value = 1637
$ws.Range("P6").Value = 105
$ws.Range("Q6").Value = 62
$ws.Range("R6").Value = -63
$ws.Range("S6").Value = -8
$ws.Range("T6").Value = 31
$ws.Range("U6").Value = 31
$ws.Range("V6").Value = 674
$ws.Range("W6").Value = 0.04
$ws.Range("X6").Value = -8.14
$ws.Range("Y6").Value = -2.21
$ws.Range("Z6").Value = -1.42
$ws.Range("AA6").Value = 60.94
$ws.Range("AB6").Value = 1455.94
$ws.Range("AC6").Value = -338
$ws.Range("AD6").Value = -38.9
$ws.Range("AE6").Value = 15221
$ws.Range("AF6").Value = 0.86
$ws.Range("AG6").Value = 191
$ws.Range("AH6").Value = 1.45
$ws.Range("AI6").Value = -55.09
$ws.Range("AJ6").Value = 11002587

# Rows 7-9: clear all financial data cells, keep A/B/C columns intact
$ws.Range("D7:AJ9").ClearContents()
